$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (data currently spans rows 1..201, columns A..H)
$lastRow = $ws.UsedRange.Rows.Count

# Step 1: capture the current marital_status column (G, the 7th column) values
# for each data row before we shift columns around.
$maritalValues = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $maritalValues[$r] = $ws.Cells.Item($r, 7).Value2
}

# Step 2: delete column F (citizenship). This shifts the old G (marital_status)
# column into F, and the old H (party) column into G - matching the target layout.
$ws.Columns("F").Delete()

# Step 3: overwrite column F (rows 2..lastRow) with the generalised marital status:
# "Never married" stays as-is, every other prior status becomes "Ever married".
for ($r = 2; $r -le $lastRow; $r++) {
    $old = $maritalValues[$r]
    if ($old -eq "Never married") {
        $ws.Cells.Item($r, 6).Value = "Never married"
    } else {
        $ws.Cells.Item($r, 6).Value = "Ever married"
    }
}
